# 202510 HL Maintain Report - append 5 new ticket rows (98-102) to the
# Report sheet, fix the wrap formatting on row 97's P/AC cells, and grow
# the print area / dimension to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# --- 1. Row 97 formatting fix: P97 / AC97 switch to the wrapped, left
#        aligned look (matches the rest of the "work content" column).
$ws.Range("P97").WrapText = $true
$ws.Range("P97").HorizontalAlignment = -4131
$ws.Range("AC97").WrapText = $true
$ws.Range("AC97").HorizontalAlignment = -4131

# --- 2. Build out rows 98-102 by cloning the existing banding (odd rows
#        use the shaded style, even rows the plain style) so the new
#        rows pick up matching fills/borders, then stamp in the data.
$ws.Range("A96:AK96").Copy($ws.Range("A98:AK98"))
$ws.Range("A97:AK97").Copy($ws.Range("A99:AK99"))
$ws.Range("A98:AK98").Copy($ws.Range("A100:AK100"))
$ws.Range("A99:AK99").Copy($ws.Range("A101:AK101"))
$ws.Range("A96:AK96").Copy($ws.Range("A102:AK102"))

# --- 3. Row 98 data ---
$ws.Range("A98").Value = 96
$ws.Range("B98").Value = "服務"
$ws.Range("C98").Value = 2025103336
$ws.Range("D98").Value = ""
$ws.Range("E98").Value = ""
$ws.Range("F98").Value = "D649"
$ws.Range("G98").Value = "三重運動公園"
$ws.Range("H98").Value = "新北市三重區"
$ws.Range("I98").Value = ""
$ws.Range("J98").Value = ""
$ws.Range("K98").Value = ""
$ws.Range("L98").Value = ""
$ws.Range("M98").Value = ""
$ws.Range("N98").Value = ""
$ws.Range("O98").Value = ""
$ws.Range("P98").Value = ""
$ws.Range("Q98").Value = "THILF0D649"
$ws.Range("R98").Value = "新北一"
$ws.Range("S98").Value = "吳宗鴻"
$ws.Range("T98").Value = 1
$ws.Range("U98").Value = "已完工"
$ws.Range("V98").Value = "2025-10-28 11:13:24"
$ws.Range("W98").Value = "2025-10-28 10:30:00"
$ws.Range("X98").Value = "2025-10-28 10:50:00"
$ws.Range("Y98").Value = ""
$ws.Range("Z98").Value = 0.3
$ws.Range("AA98").Value = ""
$ws.Range("AB98").Value = "到場處理"
$ws.Range("AC98").Value = "PMQ4+L90"
$ws.Range("AD98").Value = "O"
$ws.Range("AE98").Value = ""
$ws.Range("AF98").Value = ""
$ws.Range("AG98").Value = ""
$ws.Range("AH98").Value = ""
$ws.Range("AI98").Value = ""
$ws.Range("AJ98").Value = ""
$ws.Range("AK98").Value = "O"

# --- 4. Row 99 data ---
$ws.Range("A99").Value = 97
$ws.Range("B99").Value = "服務"
$ws.Range("C99").Value = 2025103337
$ws.Range("D99").Value = ""
$ws.Range("E99").Value = ""
$ws.Range("F99").Value = 3606
$ws.Range("G99").Value = "北縣泰富店"
$ws.Range("H99").Value = "新北市泰山區"
$ws.Range("I99").Value = ""
$ws.Range("J99").Value = ""
$ws.Range("K99").Value = ""
$ws.Range("L99").Value = ""
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = ""
$ws.Range("O99").Value = ""
$ws.Range("P99").Value = ""
$ws.Range("Q99").Value = "THILF03606"
$ws.Range("R99").Value = "新北一"
$ws.Range("S99").Value = "湯家瑋"
$ws.Range("T99").Value = 1
$ws.Range("U99").Value = "已完工"
$ws.Range("V99").Value = "2025-10-28 11:14:45"
$ws.Range("W99").Value = "2025-10-28 10:30:00"
$ws.Range("X99").Value = "2025-10-28 11:00:00"
$ws.Range("Y99").Value = ""
$ws.Range("Z99").Value = 0.5
$ws.Range("AA99").Value = ""
$ws.Range("AB99").Value = "到場處理"
$ws.Range("AC99").Value = "PMQ4"
$ws.Range("AD99").Value = "O"
$ws.Range("AE99").Value = ""
$ws.Range("AF99").Value = ""
$ws.Range("AG99").Value = ""
$ws.Range("AH99").Value = ""
$ws.Range("AI99").Value = ""
$ws.Range("AJ99").Value = ""
$ws.Range("AK99").Value = "O"

# --- 5. Row 100 data ---
$ws.Range("A100").Value = 98
$ws.Range("B100").Value = "服務"
$ws.Range("C100").Value = 2025103342
$ws.Range("D100").Value = ""
$ws.Range("E100").Value = ""
$ws.Range("F100").Value = 4656
$ws.Range("G100").Value = "泰山新民店"
$ws.Range("H100").Value = "新北市泰山區"
$ws.Range("I100").Value = ""
$ws.Range("J100").Value = ""
$ws.Range("K100").Value = ""
$ws.Range("L100").Value = ""
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = ""
$ws.Range("O100").Value = ""
$ws.Range("P100").Value = ""
$ws.Range("Q100").Value = "THILF04656"
$ws.Range("R100").Value = "新北一"
$ws.Range("S100").Value = "湯家瑋"
$ws.Range("T100").Value = 1
$ws.Range("U100").Value = "已完工"
$ws.Range("V100").Value = "2025-10-28 11:33:15"
$ws.Range("W100").Value = "2025-10-28 11:00:00"
$ws.Range("X100").Value = "2025-10-28 11:30:00"
$ws.Range("Y100").Value = ""
$ws.Range("Z100").Value = 0.5
$ws.Range("AA100").Value = ""
$ws.Range("AB100").Value = "到場處理"
$ws.Range("AC100").Value = "PMQ4+L90"
$ws.Range("AD100").Value = "O"
$ws.Range("AE100").Value = ""
$ws.Range("AF100").Value = ""
$ws.Range("AG100").Value = ""
$ws.Range("AH100").Value = ""
$ws.Range("AI100").Value = ""
$ws.Range("AJ100").Value = "O"
$ws.Range("AK100").Value = "O"

# --- 6. Row 101 data ---
$ws.Range("A101").Value = 99
$ws.Range("B101").Value = "服務"
$ws.Range("C101").Value = 2025103348
$ws.Range("D101").Value = ""
$ws.Range("E101").Value = ""
$ws.Range("F101").Value = 5377
$ws.Range("G101").Value = "三重大智店"
$ws.Range("H101").Value = "新北市三重區"
$ws.Range("I101").Value = ""
$ws.Range("J101").Value = ""
$ws.Range("K101").Value = ""
$ws.Range("L101").Value = ""
$ws.Range("M101").Value = ""
$ws.Range("N101").Value = ""
$ws.Range("O101").Value = ""
$ws.Range("P101").Value = ""
$ws.Range("Q101").Value = "THILF05377"
$ws.Range("R101").Value = "新北一"
$ws.Range("S101").Value = "吳宗鴻"
$ws.Range("T101").Value = 1
$ws.Range("U101").Value = "已完工"
$ws.Range("V101").Value = "2025-10-28 11:54:28"
$ws.Range("W101").Value = "2025-10-28 11:05:00"
$ws.Range("X101").Value = "2025-10-28 11:20:00"
$ws.Range("Y101").Value = ""
$ws.Range("Z101").Value = 0.3
$ws.Range("AA101").Value = ""
$ws.Range("AB101").Value = "到場處理"
$ws.Range("AC101").Value = "L90"
$ws.Range("AD101").Value = ""
$ws.Range("AE101").Value = ""
$ws.Range("AF101").Value = ""
$ws.Range("AG101").Value = ""
$ws.Range("AH101").Value = ""
$ws.Range("AI101").Value = ""
$ws.Range("AJ101").Value = "O"
$ws.Range("AK101").Value = "O"

# --- 7. Row 102 data ---
$ws.Range("A102").Value = 100
$ws.Range("B102").Value = "服務"
$ws.Range("C102").Value = 2025103349
$ws.Range("D102").Value = ""
$ws.Range("E102").Value = ""
$ws.Range("F102").Value = "D028"
$ws.Range("G102").Value = "三重正義北店"
$ws.Range("H102").Value = "新北市三重區"
$ws.Range("I102").Value = ""
$ws.Range("J102").Value = ""
$ws.Range("K102").Value = ""
$ws.Range("L102").Value = ""
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = ""
$ws.Range("O102").Value = ""
$ws.Range("P102").Value = ""
$ws.Range("Q102").Value = "THILF0D028"
$ws.Range("R102").Value = "新北一"
$ws.Range("S102").Value = "吳宗鴻"
$ws.Range("T102").Value = 1
$ws.Range("U102").Value = "已完工"
$ws.Range("V102").Value = "2025-10-28 12:32:58"
$ws.Range("W102").Value = "2025-10-28 12:00:00"
$ws.Range("X102").Value = "2025-10-28 12:20:00"
$ws.Range("Y102").Value = ""
$ws.Range("Z102").Value = 0.3
$ws.Range("AA102").Value = ""
$ws.Range("AB102").Value = "到場處理"
$ws.Range("AC102").Value = "PMQ4+L90"
$ws.Range("AD102").Value = "O"
$ws.Range("AE102").Value = ""
$ws.Range("AF102").Value = ""
$ws.Range("AG102").Value = ""
$ws.Range("AH102").Value = ""
$ws.Range("AI102").Value = ""
$ws.Range("AJ102").Value = ""
$ws.Range("AK102").Value = "O"

# --- 8. Grow the print area to cover the 5 new rows ---
$printArea = $wb.Names.Item("Report!Print_Area")
$printArea.RefersTo = "='Report'!`$A`$1:`$AK`$102"

# --- 9. Match the selection left behind by the edit ---
$ws.Range("AC99").Select()
